$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - bold/centered style "1", matching existing header cells
$ws.Range("AB1").Value = "Q_Responsibility"
$ws.Range("AC1").Value = "Q_Responsibilitycode"

# Data rows: AB = free-text answer, AC = numeric-looking code stored as text
# (matches existing "...code" columns, e.g. G/I/M/... which are all text cells)
$ws.Range("AB2").Value = "Overheidsinstanties en burgers zijn even verantwoordelijk voor bescherming tegen overstromingen."
$ws.Range("AB3").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC3")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB4").Value = "Overheidsinstanties zijn volledig verantwoordelijk voor bescherming tegen overstromingen."
$ws.Range("AB5").Value = "Overheidsinstanties en burgers zijn even verantwoordelijk voor bescherming tegen overstromingen."
$ws.Range("AB6").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC6")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB7").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC7")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB8").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC8")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB9").Value = "Overheidsinstanties zijn volledig verantwoordelijk voor bescherming tegen overstromingen."
$ws.Range("AB10").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC10")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB11").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC11")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB12").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC12")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB13").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC13")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB14").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC14")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB15").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC15")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB16").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC16")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB17").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC17")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB18").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC18")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB19").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC19")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB20").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC20")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB21").Value = "Overheidsinstanties en burgers zijn even verantwoordelijk voor bescherming tegen overstromingen."
$ws.Range("AB22").Value = "Overheidsinstanties en burgers zijn even verantwoordelijk voor bescherming tegen overstromingen."
$ws.Range("AB23").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC23")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
$ws.Range("AB24").Value = "Overheidsinstanties zijn verantwoordelijk en burgers deels verantwoordelijk voor bescherming tegen overstromingen"
$c = $ws.Range("AC24")
$c.NumberFormat = "@"
$c.Value = "2"
$c.Style = "Normal"
